$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after row 25 (shifts rows 26+ down by one).
# The inserted row becomes the new row 26 and inherits what used to be
# row 25's original content (A=46026, B="妙一師姐", C="聖務", D=empty).
$ws.Rows.Item(26).Insert()

# Copy original row 25 values/content down into the newly inserted row 26.
$ws.Range("A26").Value = $ws.Range("A25").Value
$ws.Range("B26").Value = $ws.Range("B25").Value
$ws.Range("C26").Value = $ws.Range("C25").Value
$ws.Range("D26").Value = ""

# Update row 24: add Note text, enable wrap + taller row.
$ws.Range("D24").Value = "課程 : 知道、求道、行道`nhttps://www.youtube.com/watch?v=Lq4ziyNg63U"
$ws.Range("D24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 32.75

# Update row 25: Shift -> 講經說法, add Note text, enable wrap + taller row.
$ws.Range("C25").Value = "講經說法"
$ws.Range("D25").Value = "課程 : 知道、求道、行道`nhttps://www.youtube.com/watch?v=Lq4ziyNg63U"
$ws.Range("D25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 32.75

# Restore selection to D21 (matches the saved file's cursor position).
$ws.Range("D21").Select()
